# Apply the "2021-04-27" model-holdings refresh to FAST_holdings.xlsx
#   - bump the confidentiality banner date (2021-04-26 -> 2021-04-27)
#   - update the Weight (col D) / Percent Change (col E) figures for rows 2-10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; unprotect (same password used in the file),
# make the edits, then restore protection.
$ws.Unprotect("D382")

# --- Confidentiality banner text (shared string used by A13) ---
$newBanner = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."
$ws.Range("A13").Value = $newBanner
# Setting a hard-wrapped (multi-line) value auto-expands the row height;
# AutoFit puts it straight back to the sheet's standard height so row 13
# ends up unchanged, exactly like the original file.
$ws.Rows.Item(13).AutoFit()

# --- Weight / Percent Change refresh (rows 2-9) ---
$updates = @(
    @{ Row = 2;  D = 0.104789703184693;    E = -0.001949317738791478 },
    @{ Row = 3;  D = 0.1094672583676208;   E = -0.003302773190592911 },
    @{ Row = 4;  D = 0.1161961288834615;   E = -0.0007840677434531873 },
    @{ Row = 5;  D = 0.1372396129181073;   E = 0.001504712124811736 },
    @{ Row = 6;  D = 0.1325083533599373;   E = -0.002358163406852465 },
    @{ Row = 7;  D = 0.1403832368831791;   E = 0.00181713848508025 },
    @{ Row = 8;  D = 0.1295490185158605;   E = -0.002837684449489286 },
    @{ Row = 9;  D = 0.1298666878871407;   E = -0.002700528951431691 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Row 10 (Total row) only has its Percent Change (col E) refreshed.
$ws.Cells.Item(10, 5).Value = -0.001226122033804011

# Restore the original sheet protection.
$ws.Protect("D382")
